$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header "Profile code" -> "Profile Name"
$ws.Range("A1").Value = "Profile Name"

# Add three new Polish (PEF.PL) process rows
$ws.Range("A29").Value = "PEF.PL Accounting Note v1"
$ws.Range("C29").Value = "cenbii-procid-ubl"
$ws.Range("D29").Value = "urn:fdc:www.efaktura.gov.pl:ver1.0:account_corr:ver1.0"
$ws.Range("E29").Value = "6"
$ws.Range("F29").Value = $false

$ws.Range("A30").Value = "PEF.PL Correcting Invoice v1"
$ws.Range("C30").Value = "cenbii-procid-ubl"
$ws.Range("D30").Value = "urn:fdc:www.efaktura.gov.pl:ver1.0:corr_inv:ver1.0"
$ws.Range("E30").Value = "6"
$ws.Range("F30").Value = $false

$ws.Range("A31").Value = "PEF.PL Receipt Advice v1"
$ws.Range("C31").Value = "cenbii-procid-ubl"
$ws.Range("D31").Value = "urn:fdc:www.efaktura.gov.pl:ver1.0:receipt_advice:ver1.0"
$ws.Range("E31").Value = "6"
$ws.Range("F31").Value = $false

# Update selection to reflect the new last-edited cells
$ws.Range("E31:F31").Select()
